$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Remove the first worker's (FERNANDO ISSAIAS MORALES ESPINOSA)
#    seven data rows entirely. Excel shifts everything below
#    (the second worker's rows, the blank spacer rows and the
#    signature block) up by 7 rows automatically, which also moves
#    the merged cells and the drawing anchor along with it.
# ------------------------------------------------------------------
$ws.Rows("16:22").Delete()

# ------------------------------------------------------------------
# 2. The remaining worker (MARIA BERNARDA NOVA MARRUGO) now occupies
#    rows 16-22. Re-order her periods from descending (2408 -> 2402)
#    to ascending (2402 -> 2408) and refresh the "Valor Mora" /
#    "Salario Basico" figures to the new period-2 amounts.
# ------------------------------------------------------------------
$periods = @("2402", "2403", "2404", "2405", "2406", "2407", "2408")
$valorMora = @(320000, 320000, 320000, 320000, 320000, 320000, 224000)
$salarioBasico = @(8000000, 8000000, 8000000, 8000000, 8000000, 8000000, 8000000)

for ($i = 0; $i -lt 7; $i++) {
    $r = 16 + $i
    $ws.Cells.Item($r, 5).Value = $periods[$i]
    $ws.Cells.Item($r, 6).Value = $valorMora[$i]
    $ws.Cells.Item($r, 7).Value = $salarioBasico[$i]
}

# ------------------------------------------------------------------
# 3. Update the summary figures above the table: total overdue value
#    and the worker count (only one worker remains now).
# ------------------------------------------------------------------
$ws.Range("E11").Value = 2144000
$ws.Range("C13").Value = 1

# ------------------------------------------------------------------
# 4. Let Excel recompute the "best fit" width of column D now that
#    the longer worker name (FERNANDO ISSAIAS MORALES ESPINOSA) is
#    gone and the remaining name (MARIA BERNARDA NOVA MARRUGO) is
#    shorter, so the column shrinks to fit it.
# ------------------------------------------------------------------
$ws.Columns("D").AutoFit()
$ws.Columns("D").ColumnWidth = 31.8
